{"js": "/*\n * Applies the addition/subtraction worksheet update:\n *  - Updates the date heading paragraph.\n *  - Updates every arithmetic expression in the 20x5 table, cell by cell,\n *    in document order (paragraphs collection includes the heading\n *    paragraph followed by every table-cell paragraph, in reading order).\n *\n * Both the \"before\" and \"after\" text for every paragraph are listed\n * explicitly (taken from the canonical OOXML diff) so the replacement is\n * applied positionally and safely handles the few duplicate expressions\n * (e.g. \"24+66=90\" and \"33+17=50\" each occur twice, with different\n * replacements at each occurrence).\n */\nconst oldTexts = [\"2024-12-17 Tuesday\", \"88-59=29\", \"78-32=46\", \"82-65=17\", \"4+35=39\", \"15+23=38\", \"1+62=63\", \"71+10=81\", \"98-7=91\", \"87+5=92\", \"80+18=98\", \"44+8=52\", \"41-13=28\", \"19-7=12\", \"4+52=56\", \"37+17=54\", \"75-21=54\", \"23+16=39\", \"26+38=64\", \"15+32=47\", \"24+66=90\", \"87-61=26\", \"67-22=45\", \"99-51=48\", \"15+41=56\", \"62-45=17\", \"55+40=95\", \"1+80=81\", \"32-25=7\", \"19+25=44\", \"14-4=10\", \"33+17=50\", \"0+41=41\", \"33+17=50\", \"15+29=44\", \"55-48=7\", \"76-75=1\", \"67-63=4\", \"14+50=64\", \"86+4=90\", \"80-38=42\", \"66-57=9\", \"33+16=49\", \"87-60=27\", \"86-28=58\", \"24+66=90\", \"27+4=31\", \"14+51=65\", \"19+11=30\", \"96-11=85\", \"57-23=34\", \"25+27=52\", \"86-4=82\", \"84-3=81\", \"72+27=99\", \"81-14=67\", \"7+37=44\", \"50-18=32\", \"44+32=76\", \"45-37=8\", \"74-52=22\", \"87-9=78\", \"54-49=5\", \"80-30=50\", \"13+52=65\", \"4+36=40\", \"16+69=85\", \"59+14=73\", \"56+21=77\", \"9+9=18\", \"12+67=79\", \"53-24=29\", \"51-8=43\", \"49-45=4\", \"83-24=59\", \"78-28=50\", \"9-4=5\", \"99-52=47\", \"77-45=32\", \"96-9=87\", \"41+36=77\", \"36+6=42\", \"64-50=14\", \"31+66=97\", \"68-45=23\", \"30-1=29\", \"79-20=59\", \"57-20=37\", \"40-17=23\", \"94-40=54\", \"93-44=49\", \"47+27=74\", \"72-6=66\", \"99-93=6\", \"13+31=44\", \"3+61=64\", \"67-17=50\", \"94-54=40\", \"34+46=80\", \"81+13=94\", \"40+6=46\"];\nconst newTexts = [\"2024-12-18 Wednesday\", \"83-59=24\", \"5+38=43\", \"89-9=80\", \"7+23=30\", \"18+5=23\", \"81-21=60\", \"61+35=96\", \"27-11=16\", \"2+64=66\", \"2+90=92\", \"38-4=34\", \"92-63=29\", \"56-9=47\", \"54-34=20\", \"49-34=15\", \"40-11=29\", \"81-11=70\", \"15+15=30\", \"46-45=1\", \"58-11=47\", \"24-2=22\", \"4+29=33\", \"71-63=8\", \"13+76=89\", \"71+22=93\", \"93-13=80\", \"69-55=14\", \"82-57=25\", \"10+16=26\", \"63-54=9\", \"28-18=10\", \"80-12=68\", \"4+48=52\", \"33+21=54\", \"76-52=24\", \"71+6=77\", \"48-25=23\", \"86-3=83\", \"70-45=25\", \"92-83=9\", \"68-4=64\", \"74-70=4\", \"64+27=91\", \"35-25=10\", \"88-84=4\", \"43+4=47\", \"72-64=8\", \"94-71=23\", \"23-10=13\", \"64-41=23\", \"30+18=48\", \"11+6=17\", \"90+6=96\", \"34+28=62\", \"65-4=61\", \"52-4=48\", \"28+1=29\", \"69+9=78\", \"6+91=97\", \"10+57=67\", \"27+52=79\", \"77-49=28\", \"66-44=22\", \"36+37=73\", \"91-47=44\", \"78+19=97\", \"74-48=26\", \"92-59=33\", \"51+7=58\", \"39-26=13\", \"24-1=23\", \"43-14=29\", \"42+51=93\", \"35+14=49\", \"8+42=50\", \"10+48=58\", \"68+21=89\", \"62+27=89\", \"84-37=47\", \"24+73=97\", \"80-60=20\", \"40-34=6\", \"78-63=15\", \"27-13=14\", \"13+65=78\", \"14+69=83\", \"12+66=78\", \"10+76=86\", \"68-62=6\", \"95-36=59\", \"18+38=56\", \"46+37=83\", \"39+14=53\", \"60-54=6\", \"95-25=70\", \"91-30=61\", \"48-17=31\", \"18+69=87\", \"9+43=52\", \"88-62=26\"];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== oldTexts.length) {\n  throw new Error(\n    `Expected ${oldTexts.length} paragraphs (1 heading + 100 table cells), found ${items.length}`\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const para = items[i];\n  const expected = oldTexts[i];\n  const actual = para.text;\n  // Guard against drift between this script's assumptions and the live\n  // document content; only touch paragraphs whose text still matches what\n  // the diff recorded as the \"before\" value.\n  if (actual === expected) {\n    if (newTexts[i] !== expected) {\n      para.insertText(newTexts[i], \"Replace\");\n    }\n  } else if (actual !== newTexts[i]) {\n    // Text differs from both expected old and new values -- surface this\n    // rather than silently skipping it.\n    throw new Error(\n      `Paragraph ${i}: expected \"${expected}\" but found \"${actual}\"`\n    );\n  }\n}\n\nawait context.sync();\n", "ps1": "# Applies the addition/subtraction worksheet update:\n#  - Updates the date heading paragraph.\n#  - Updates every arithmetic expression in the 20x5 table, cell by cell,\n#    addressed positionally via Table.Cell(row, col) so the handful of\n#    duplicate expressions (e.g. \"24+66=90\" and \"33+17=50\" each appear twice,\n#    with different replacements at each occurrence) are handled correctly.\n\n$d = $word.ActiveDocument\n\n$oldDate = \"2024-12-17 Tuesday\"\n$newDate = \"2024-12-18 Wednesday\"\n\n# 20 rows x 5 columns of arithmetic expressions, in reading order\n# (left-to-right, top-to-bottom), taken from the canonical OOXML diff.\n$oldCells = @(\n    @(\"88-59=29\", \"78-32=46\", \"82-65=17\", \"4+35=39\", \"15+23=38\"),\n    @(\"1+62=63\", \"71+10=81\", \"98-7=91\", \"87+5=92\", \"80+18=98\"),\n    @(\"44+8=52\", \"41-13=28\", \"19-7=12\", \"4+52=56\", \"37+17=54\"),\n    @(\"75-21=54\", \"23+16=39\", \"26+38=64\", \"15+32=47\", \"24+66=90\"),\n    @(\"87-61=26\", \"67-22=45\", \"99-51=48\", \"15+41=56\", \"62-45=17\"),\n    @(\"55+40=95\", \"1+80=81\", \"32-25=7\", \"19+25=44\", \"14-4=10\"),\n    @(\"33+17=50\", \"0+41=41\", \"33+17=50\", \"15+29=44\", \"55-48=7\"),\n    @(\"76-75=1\", \"67-63=4\", \"14+50=64\", \"86+4=90\", \"80-38=42\"),\n    @(\"66-57=9\", \"33+16=49\", \"87-60=27\", \"86-28=58\", \"24+66=90\"),\n    @(\"27+4=31\", \"14+51=65\", \"19+11=30\", \"96-11=85\", \"57-23=34\"),\n    @(\"25+27=52\", \"86-4=82\", \"84-3=81\", \"72+27=99\", \"81-14=67\"),\n    @(\"7+37=44\", \"50-18=32\", \"44+32=76\", \"45-37=8\", \"74-52=22\"),\n    @(\"87-9=78\", \"54-49=5\", \"80-30=50\", \"13+52=65\", \"4+36=40\"),\n    @(\"16+69=85\", \"59+14=73\", \"56+21=77\", \"9+9=18\", \"12+67=79\"),\n    @(\"53-24=29\", \"51-8=43\", \"49-45=4\", \"83-24=59\", \"78-28=50\"),\n    @(\"9-4=5\", \"99-52=47\", \"77-45=32\", \"96-9=87\", \"41+36=77\"),\n    @(\"36+6=42\", \"64-50=14\", \"31+66=97\", \"68-45=23\", \"30-1=29\"),\n    @(\"79-20=59\", \"57-20=37\", \"40-17=23\", \"94-40=54\", \"93-44=49\"),\n    @(\"47+27=74\", \"72-6=66\", \"99-93=6\", \"13+31=44\", \"3+61=64\"),\n    @(\"67-17=50\", \"94-54=40\", \"34+46=80\", \"81+13=94\", \"40+6=46\")\n)\n$newCells = @(\n    @(\"83-59=24\", \"5+38=43\", \"89-9=80\", \"7+23=30\", \"18+5=23\"),\n    @(\"81-21=60\", \"61+35=96\", \"27-11=16\", \"2+64=66\", \"2+90=92\"),\n    @(\"38-4=34\", \"92-63=29\", \"56-9=47\", \"54-34=20\", \"49-34=15\"),\n    @(\"40-11=29\", \"81-11=70\", \"15+15=30\", \"46-45=1\", \"58-11=47\"),\n    @(\"24-2=22\", \"4+29=33\", \"71-63=8\", \"13+76=89\", \"71+22=93\"),\n    @(\"93-13=80\", \"69-55=14\", \"82-57=25\", \"10+16=26\", \"63-54=9\"),\n    @(\"28-18=10\", \"80-12=68\", \"4+48=52\", \"33+21=54\", \"76-52=24\"),\n    @(\"71+6=77\", \"48-25=23\", \"86-3=83\", \"70-45=25\", \"92-83=9\"),\n    @(\"68-4=64\", \"74-70=4\", \"64+27=91\", \"35-25=10\", \"88-84=4\"),\n    @(\"43+4=47\", \"72-64=8\", \"94-71=23\", \"23-10=13\", \"64-41=23\"),\n    @(\"30+18=48\", \"11+6=17\", \"90+6=96\", \"34+28=62\", \"65-4=61\"),\n    @(\"52-4=48\", \"28+1=29\", \"69+9=78\", \"6+91=97\", \"10+57=67\"),\n    @(\"27+52=79\", \"77-49=28\", \"66-44=22\", \"36+37=73\", \"91-47=44\"),\n    @(\"78+19=97\", \"74-48=26\", \"92-59=33\", \"51+7=58\", \"39-26=13\"),\n    @(\"24-1=23\", \"43-14=29\", \"42+51=93\", \"35+14=49\", \"8+42=50\"),\n    @(\"10+48=58\", \"68+21=89\", \"62+27=89\", \"84-37=47\", \"24+73=97\"),\n    @(\"80-60=20\", \"40-34=6\", \"78-63=15\", \"27-13=14\", \"13+65=78\"),\n    @(\"14+69=83\", \"12+66=78\", \"10+76=86\", \"68-62=6\", \"95-36=59\"),\n    @(\"18+38=56\", \"46+37=83\", \"39+14=53\", \"60-54=6\", \"95-25=70\"),\n    @(\"91-30=61\", \"48-17=31\", \"18+69=87\", \"9+43=52\", \"88-62=26\")\n)\n\n# --- Update the date heading ---------------------------------------------\n$dateParagraph = $d.Paragraphs.Item(1)\n$dateText = $dateParagraph.Range.Text.TrimEnd([char]13, [char]7)\nif ($dateText -eq $oldDate) {\n    $dateParagraph.Range.Text = $newDate\n} elseif ($dateText -ne $newDate) {\n    throw \"Date paragraph: expected `\"$oldDate`\" but found `\"$dateText`\"\"\n}\n\n# --- Update every cell in the practice-problems table ---------------------\n$table = $d.Tables.Item(1)\nfor ($r = 1; $r -le 20; $r++) {\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $table.Cell($r, $c)\n        $expected = $oldCells[$r - 1][$c - 1]\n        $replacement = $newCells[$r - 1][$c - 1]\n        $actual = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($actual -eq $expected) {\n            if ($replacement -ne $expected) {\n                $cell.Range.Text = $replacement\n            }\n        } elseif ($actual -ne $replacement) {\n            throw \"Cell ($r,$c): expected `\"$expected`\" but found `\"$actual`\"\"\n        }\n    }\n}\n"}
